# Reorders the "Observation"/"Contact" data blocks in the dataDefinition
# sheet. Each block is a small run of rows (one "header" row carrying
# Kategori/Objekt/Beskrivelse in columns A/B/C, followed by 1-4 "detail"
# rows carrying Variable/Type/.../Required in columns D-H).
#
# Blocks, identified by their Objekt name (col B of the header row) and
# their *original* row span (inclusive):
#   SystolicBloodPressure  2-4
#   OxygenDemand           5-7
#   AVPUScale              8-10
#   Triage                 11-13
#   RespiratoryRate        14-16
#   Temperature            17-19
#   GlasgowComaScale       20-22
#   PulseRate              23-25
#   OxygenSaturation       26-28
#   PainEvaluation         29-33
#   NoteType (Contact)     34-36
#
# New order (top to bottom):
#   RespiratoryRate, AVPUScale, PulseRate, OxygenDemand, OxygenSaturation,
#   Triage, GlasgowComaScale, Temperature, PainEvaluation,
#   SystolicBloodPressure, NoteType
#
# Strategy: copy cell-by-cell with Range.Copy (never Value assignment) so
# cell types are preserved exactly (notably the literal text "True" in
# column H, which must stay text and not be coerced to a boolean), and
# only touch cells that actually hold data (so we don't litter the sheet
# with empty placeholder cells for the blank columns every block leaves
# untouched). Blocks are first copied to a scratch area far below the
# used range, the original rows are cleared, then the scratch copies are
# written back in the new order. Finally the scratch area is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numCols = 8

$blockStart = @{
    "SystolicBloodPressure" = 2
    "OxygenDemand" = 5
    "AVPUScale" = 8
    "Triage" = 11
    "RespiratoryRate" = 14
    "Temperature" = 17
    "GlasgowComaScale" = 20
    "PulseRate" = 23
    "OxygenSaturation" = 26
    "PainEvaluation" = 29
    "NoteType" = 34
}
$blockEnd = @{
    "SystolicBloodPressure" = 4
    "OxygenDemand" = 7
    "AVPUScale" = 10
    "Triage" = 13
    "RespiratoryRate" = 16
    "Temperature" = 19
    "GlasgowComaScale" = 22
    "PulseRate" = 25
    "OxygenSaturation" = 28
    "PainEvaluation" = 33
    "NoteType" = 36
}

# Target order (top to bottom) after the edit.
$newOrder = @("RespiratoryRate","AVPUScale","PulseRate","OxygenDemand","OxygenSaturation","Triage","GlasgowComaScale","Temperature","PainEvaluation","SystolicBloodPressure","NoteType")

$firstDataRow = 2
$lastDataRow = 36
$scratchStart = 100

function Copy-BlockCellwise($ws, $srcStartRow, $numRows, $dstStartRow, $numCols) {
    for ($i = 0; $i -lt $numRows; $i++) {
        for ($c = 1; $c -le $numCols; $c++) {
            $srcCell = $ws.Cells.Item($srcStartRow + $i, $c)
            if ($srcCell.Value2 -ne $null) {
                $dstCell = $ws.Cells.Item($dstStartRow + $i, $c)
                $srcCell.Copy($dstCell)
            }
        }
    }
}

# --- Step 1: copy every block (in new-order sequence) down into a scratch
#     area, so we have an independent copy of each block's cells/types
#     before we touch the live rows.
$scratchRow = @{}
$cursor = $scratchStart
foreach ($name in $newOrder) {
    $s = $blockStart[$name]
    $e = $blockEnd[$name]
    $len = $e - $s + 1
    Copy-BlockCellwise $ws $s $len $cursor $numCols
    $scratchRow[$name] = $cursor
    $cursor = $cursor + $len
}
$scratchEnd = $cursor - 1

# --- Step 2: clear the original block rows completely (ClearContents
#     removes the cell nodes entirely rather than leaving them as blanks).
$ws.Range("A" + $firstDataRow + ":H" + $lastDataRow).ClearContents()

# --- Step 3: copy each block from the scratch area back into its new
#     home, walking down from row 2 in the new order.
$writeRow = $firstDataRow
foreach ($name in $newOrder) {
    $s = $blockStart[$name]
    $e = $blockEnd[$name]
    $len = $e - $s + 1
    $sRow = $scratchRow[$name]
    Copy-BlockCellwise $ws $sRow $len $writeRow $numCols
    $writeRow = $writeRow + $len
}

# --- Step 4: clear the scratch area, leaving the sheet's used range back
#     to its original extent.
$ws.Range("A" + $scratchStart + ":H" + $scratchEnd).ClearContents()
